$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> (Price, Volume) updates. Empty string means "no change".
$updates = @(
    @{ Row = 2;  D = "26.957.22";  E = "  -1.81%  " },
    @{ Row = 3;  D = "1.565.88";   E = "  -0.30%  " },
    @{ Row = 4;  D = "";           E = "  +0.26%  " },
    @{ Row = 5;  D = "206.29";     E = "  -0.56%  " },
    @{ Row = 6;  D = "0.488";      E = "  -1.69%  " },
    @{ Row = 7;  D = "";           E = "  +0.19%  " },
    @{ Row = 8;  D = "22.12";      E = "  +0.13%  " },
    @{ Row = 9;  D = "";           E = "  -0.19%  " },
    @{ Row = 10; D = "";           E = "  -0.66%  " },
    @{ Row = 11; D = "";           E = "  +0.26%  " },
    @{ Row = 12; D = "1.788.28";   E = "  -0.39%  " },
    @{ Row = 13; D = "1.558.83";   E = "  -1.02%  " },
    @{ Row = 14; D = "3.76";       E = "  -1.51%  " },
    @{ Row = 15; D = "0.515";      E = "  -1.03%  " },
    @{ Row = 16; D = "26.959.45";  E = "  -1.76%  " },
    @{ Row = 17; D = "61.73";      E = "  -2.24%  " },
    @{ Row = 18; D = "214.51";     E = "  +0.27%  " },
    @{ Row = 19; D = "7.36";       E = "  +1.41%  " },
    @{ Row = 20; D = "";           E = "  -1.30%  " },
    @{ Row = 21; D = "";           E = "  +0.32%  " },
    @{ Row = 22; D = "";           E = "  -0.43%  " },
    @{ Row = 23; D = "";           E = "  -3.17%  " },
    @{ Row = 24; D = "";           E = "  -0.02%  " },
    @{ Row = 25; D = "152.32";     E = "  -0.34%  " },
    @{ Row = 26; D = "6.67";       E = "  -2.47%  " },
    @{ Row = 27; D = "14.87";      E = "  -1.08%  " },
    @{ Row = 28; D = "1.00";       E = "  +0.21%  " },
    @{ Row = 29; D = "";           E = "  -1.64%  " },
    @{ Row = 30; D = "";           E = "  -3.15%  " },
    @{ Row = 31; D = "";           E = "  -2.29%  " },
    @{ Row = 32; D = "";           E = "  -1.54%  " },
    @{ Row = 33; D = "1.388.20";   E = "  +2.10%  " },
    @{ Row = 34; D = "";           E = "  -1.01%  " },
    @{ Row = 35; D = "";           E = "  +0.82%  " },
    @{ Row = 36; D = "";           E = "  -0.57%  " },
    @{ Row = 37; D = "0.943";      E = "  -2.78%  " },
    @{ Row = 38; D = "";           E = "  -2.25%  " },
    @{ Row = 39; D = "0.809";      E = "  -1.64%  " },
    @{ Row = 40; D = "0.512";      E = "  -3.79%  " },
    @{ Row = 41; D = "";           E = "  +0.27%  " },
    @{ Row = 42; D = "1.01";       E = "  +3.56%  " },
    @{ Row = 43; D = "";           E = "  +2.39%  " },
    @{ Row = 44; D = "";           E = "  +0.69%  " },
    @{ Row = 45; D = "2.19";       E = "  +1.23%  " },
    @{ Row = 46; D = "63.50";      E = "  -1.05%  " },
    @{ Row = 47; D = "1.701.72";   E = "" },
    @{ Row = 48; D = "85.35";      E = "  +0.04%  " },
    @{ Row = 49; D = "0.0$([char]8327)0971"; E = "  -2.59%  " },
    @{ Row = 50; D = "0.0495";     E = "  +0.18%  " },
    @{ Row = 51; D = "0.0949";     E = "  -0.55%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.D -ne "") {
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.E -ne "") {
        $cell = $ws.Range("E$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
